$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row labels that changed (categories replaced due to full dataset availability)
$ws.Range("A4").Value = "Breast Imaging (br)"
$ws.Range("A5").Value = "Geritourinary Radiology (gu)"
$ws.Range("A6").Value = "Ultrasound (us)"
$ws.Range("A7").Value = "Chest Radiology (ch)"
$ws.Range("A8").Value = "Interventional Radiology (ir)"

# Update numeric values for rows 2-8 (Mean, SD, Sample size)
$ws.Range("B2").Value = 0.003011139392040991
$ws.Range("C2").Value = 0.02107797574428694
$ws.Range("D2").Value = 100

$ws.Range("B3").Value = 0.003958107933239666
$ws.Range("C3").Value = 0.03857886401984387
$ws.Range("D3").Value = 96

$ws.Range("B4").Value = 0.216285536686771
$ws.Range("C4").Value = 0.1302537377122266
$ws.Range("D4").Value = 98

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 99

$ws.Range("B6").Value = 0.01136252133392825
$ws.Range("C6").Value = 0.04853106736233735
$ws.Range("D6").Value = 96

$ws.Range("B7").Value = 0.02907852821912623
$ws.Range("C7").Value = 0.05366948107284659
$ws.Range("D7").Value = 99

$ws.Range("B8").Value = 0.01112270901174773
$ws.Range("C8").Value = 0.05560858481602785
$ws.Range("D8").Value = 96

# Remove rows 9-11, which no longer exist in the updated dataset
$ws.Range("A9:D11").Delete()
